# Update simulation results (pl_mw) for the 380 kV case: columns B:O, rows 2:25
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (columns B..O), in row order starting at row 2
$rows = @(
  @(0.359424906531558,0.06656884671737373,0.07904258385524088,0.4182501562996634,0,0.4518826569307635,0.5949968088760471,0,0,0.3497777697612605,0,0,1.210060478625749,2.048610096614539),
  @(0.317904575707729,0.06105512680248637,0.0716323076375005,0.3649548114872374,0,0.4495524980890764,0.5977938548666799,0,0,0.3054756823095772,0,0,1.224854705531824,2.04924439589513),
  @(0.2924272579348326,0.05764498115483718,0.06711559542235079,0.3323125384665104,0,0.4484680340317126,0.5998015819725495,0,0,0.2782250966426147,0,0,1.234398669839506,2.050965871826449),
  @(0.2820495452264424,0.05624914768324629,0.06528336908314714,0.3190294596897587,0,0.4481129601383103,0.6006927398412785,0,0,0.2671083427489407,0,0,1.238403581576526,2.052001839548367),
  @(0.2803266174429382,0.05601699903039048,0.0649796353125538,0.3168249123538089,0,0.4480592398341301,0.6008451241508439,0,0,0.2652617057052566,0,0,1.239075578385604,2.052194046931561),
  @(0.2922872817239295,0.05762618139168296,0.06709085146442817,0.332133323649856,0,0.4484628940011248,0.5998133049133116,0,0,0.2780752199615222,0,0,1.234452213229348,2.050978489726049),
  @(0.3451055298221206,0.06467286445048614,0.07648063112706893,0.3998561154256493,0,0.4510072186724727,0.5959009637396235,0,0,0.3345127808413793,0,0,1.215065965681181,2.048552039981388),
  @(0.4487994091658152,0.07829431017587751,0.09515792090945752,0.5333802816029447,0,0.4587548462386906,0.5905334960002051,0,0,0.4447869902886339,0,0,1.18070378249717,2.05438793540506),
  @(0.5250462244110281,0.08818129690928345,0.1090429298189832,0.6320321262224553,0,0.466144870125234,0.5879969260191871,0,0,0.525554309126278,0,0,1.157687931628438,2.065173260329061),
  @(0.5597451756271994,0.09265288179744857,0.1153954165748274,0.6770561027765041,0,0.4698790738439556,0.5871489499851918,0,0,0.5622418326533136,0,0,1.147702184904412,2.071499467508346),
  @(0.5728864733998194,0.0943423819036866,0.1178061383850491,0.6941284512734285,0,0.4713469427584585,0.5868718679554377,0,0,0.5761264129855874,0,0,1.143990561707599,2.074099859309086),
  @(0.5700561974939546,0.09397868719806013,0.1172867165971354,0.6904505776565912,0,0.4710284138573115,0.5869295837710524,0,0,0.5731364912514323,0,0,1.144786821469196,2.073530700321015),
  @(0.5608262893939298,0.09279195424549869,0.115593644563674,0.678460189851279,0,0.4699987564971195,0.5871252716775359,0,0,0.5633842920200323,0,0,1.147395428569306,2.071709295398676),
  @(0.5551728940421583,0.09206455086749088,0.1145572620656452,0.6711187422276481,0,0.4693750763199489,0.5872508710346125,0,0,0.5574097091808596,0,0,1.149002364582373,2.070620320822826),
  @(0.5227788082857785,0.08788854014147773,0.1086285049565419,0.6290927966956019,0,0.4659083470213403,0.5880585023154481,0,0,0.5231555708439828,0,0,1.158350279089629,2.064788452365576),
  @(0.5029093309368022,0.08531998713245059,0.1050006418642653,0.6033499929091164,0,0.4638771905716084,0.5886323379426415,0,0,0.5021276440498923,0,0,1.164209044323753,2.061574898130715),
  @(0.4914822669456385,0.08384017365676755,0.1029173927158098,0.5885571276782429,0,0.4627439643904978,0.5889911856164076,0,0,0.490027885907665,0,0,1.167624426176221,2.059860158390023),
  @(0.4876135019271999,0.08333871488579803,0.1022126252666169,0.5835508233430176,0,0.4623662842489296,0.589117629324079,0,0,0.485930261430525,0,0,1.168788639637721,2.059302506281625),
  @(0.5050243358611795,0.08559366785284794,0.105386482173742,0.6060889270657128,0,0.4640897821910244,0.5885682719684837,0,0,0.504366628635978,0,0,1.163580651027035,2.06190315371714),
  @(0.5635372983547029,0.09314062975550996,0.1160908006466173,0.6819814269385631,0,0.4702997295216989,0.587066598188045,0,0,0.5662489750103248,0,0,1.146627323098611,2.072238723211171),
  @(0.6017877798340976,0.09805088185791533,0.1231168328913554,0.7317149135857619,0,0.4746720338893056,0.5863418146724229,0,0,0.6066447363129726,0,0,1.135954031370812,2.080187600773087),
  @(0.5813721106760283,0.09543222713433863,0.1193641567093096,0.7051584682097456,0,0.4723096623785352,0.5867051496616824,0,0,0.5850892864551156,0,0,1.141613314381596,2.075835682091707),
  @(0.5040681540537832,0.08546994651025841,0.1052120361138691,0.6048506332350314,0,0.4639935620307796,0.588597146041721,0,0,0.5033544164519697,0,0,1.163864600925862,2.061754335811287),
  @(0.4207357945730621,0.07463050427749351,0.09007676381445151,0.4971697679488756,0,0.4563619330271109,0.5917386144806045,0,0,0.4149984622792999,0,0,1.189608226904468,2.071709295398676)
)

$data = New-Object "object[,]" 24,14
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $rows[$i]
  for ($j = 0; $j -lt $r.Count; $j++) {
    $data[$i,$j] = $r[$j]
  }
}

$ws.Range("B2:O25").Value = $data

Write-Host "Updated B2:O25 with 380 kV case results"
